$wb = $excel.ActiveWorkbook

# --- Portfolio2: swap ticker SPYW.DE -> VYM in B4 ---
$ws2 = $wb.Worksheets.Item("Portfolio2")
$ws2.Range("B4").Value = "VYM"
$ws2.Activate()
$ws2.Range("B5").Select()

# --- Portfolio3: remove the ISAG.L holding row (row 8), shifting rows up ---
$ws3 = $wb.Worksheets.Item("Portfolio3")
$ws3.Rows.Item(8).Delete()
$ws3.Activate()
$ws3.Range("E12").Select()
